$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.7
$ws.Range("H2").Value = 2.6
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 4.75
$ws.Range("L2").Value = 3.5
$ws.Range("Y2").Value = 17
$ws.Range("AH2").Value = 9.5
$ws.Range("BA2").Value = 126

# Row 3
$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 2.08

# Row 4
$ws.Range("J4").Value = 2.3
$ws.Range("AM4").Value = 501

# Row 5
$ws.Range("AL5").Value = 34
$ws.Range("AY5").Value = 26
$ws.Range("AZ5").Value = 51

# Row 7
$ws.Range("G7").Value = 9.5
$ws.Range("I7").Value = 1.36
$ws.Range("L7").Value = 1.91
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.7
$ws.Range("AW7").Value = 3
